$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# This edit:
#   1. Removes the "Meta description: ..." paragraph that sits right under
#      the H1 title.
#   2. Splits the closing "Prompt: Create a feature image..." paragraph into
#      two paragraphs: a new bold "Play Chunky Fruits for Free - ..." title
#      paragraph, followed by an (still italic) paragraph whose text becomes
#      the former meta-description copy.
# ---------------------------------------------------------------------------

# --- Locate the source ("Meta description") and target ("Prompt: ...")
#     paragraphs dynamically, by content, rather than by fixed index. ---
$metaLabel = "Meta description:"
$promptLabel = "Prompt: Create a feature image"

$metaParaIndex = -1
$promptParaIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($metaParaIndex -eq -1 -and $t.StartsWith($metaLabel)) {
        $metaParaIndex = $i
    }
    if ($t.StartsWith($promptLabel)) {
        $promptParaIndex = $i
    }
}
if ($metaParaIndex -eq -1 -or $promptParaIndex -eq -1) {
    throw "Could not locate required paragraphs (meta=$metaParaIndex, prompt=$promptParaIndex)"
}

# --- STEP 1: Duplicate the "Meta description" paragraph's exact run
#     structure (it already has the desired <w:r/> + bold-run shape) to just
#     before the closing "Prompt" paragraph. Round-tripping through
#     WordOpenXML/InsertXML (rather than Range.Text/InsertParagraphBefore)
#     preserves the leading empty run that a fresh paragraph created through
#     the Range API would not carry. ---
$metaPara = $d.Paragraphs.Item($metaParaIndex)
$srcXml = $metaPara.Range.WordOpenXML
# Strip the synthetic rsid/paraId attributes WordOpenXML invents for the
# fragment so they don't leak into the destination paragraph.
$cleanXml = $srcXml -replace ' w14:paraId="[0-9A-Fa-f]*"', ''
$cleanXml = $cleanXml -replace ' w14:textId="[0-9A-Fa-f]*"', ''
$cleanXml = $cleanXml -replace ' w:rsidR="[0-9A-Fa-f]*"', ''
$cleanXml = $cleanXml -replace ' w:rsidRDefault="[0-9A-Fa-f]*"', ''

$countBeforeInsert = $d.Paragraphs.Count
$promptPara = $d.Paragraphs.Item($promptParaIndex)
$insertPoint = $d.Range($promptPara.Range.Start, $promptPara.Range.Start)
$insertPoint.InsertXML($cleanXml)
$countAfterInsert = $d.Paragraphs.Count
$numInserted = $countAfterInsert - $countBeforeInsert

# InsertXML of a single-paragraph range also carries along that range's
# trailing (end-of-story) paragraph mark, which shows up as an extra blank
# paragraph immediately after the duplicated one; remove it.
if ($numInserted -gt 1) {
    $strayIndex = $promptParaIndex + 1
    $d.Paragraphs.Item($strayIndex).Range.Delete()
}

# --- STEP 2: Trim the newly duplicated paragraph down to just its bold run,
#     retexted as the new page title. ---
$newPara = $d.Paragraphs.Item($promptParaIndex)
$pStart = $newPara.Range.Start
$pEnd = $newPara.Range.End
$labelLen = "Meta description".Length
$trailingRange = $d.Range($pStart + $labelLen, $pEnd - 1)
if ($trailingRange.Start -lt $trailingRange.End) {
    $trailingRange.Delete()
}
$boldRange = $d.Range($pStart, $pStart + $labelLen)
$boldRange.Text = "Play Chunky Fruits for Free - Exciting Wild Symbols and Modern Features"

# --- STEP 3: Remove the original "Meta description" paragraph. ---
$metaPara = $d.Paragraphs.Item($metaParaIndex)
$metaPara.Range.Delete()

# --- STEP 4: Update the final "Prompt" paragraph's text, keeping its
#     existing italic run/formatting intact. ---
$finalPromptIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text.StartsWith($promptLabel)) {
        $finalPromptIndex = $i
        break
    }
}
if ($finalPromptIndex -eq -1) {
    throw "Could not locate the Prompt paragraph after edits"
}
$promptPara = $d.Paragraphs.Item($finalPromptIndex)
$textRange = $d.Range($promptPara.Range.Start, $promptPara.Range.End - 1)
$textRange.Text = "Read our review of Chunky Fruits, a slot game with a polygonal fruit theme, expanding Wild symbols, and a mix of nostalgia and modern features. Play for free now."

Write-Host "Done. Paragraph count:" $d.Paragraphs.Count
